$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.883656666666667
$ws.Range("H2").Value = 5.650970000000001
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.495057333333333
$ws.Range("N2").Value = 16.485172
$ws.Range("O2").Value = 0.8161989011161211
$ws.Range("P2").Value = 0.8403205285996808
$ws.Range("Q2").Value = 10.35080137964889
$ws.Range("R2").Value = 93.15721241684001
$ws.Range("S2").Value = 0.8161989011161211
$ws.Range("T2").Value = 0.8403205285996808

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.883656666666667
$ws.Range("H3").Value = 5.650970000000001
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.657666
$ws.Range("N3").Value = 1.972998
$ws.Range("O3").Value = 0.09768528951377062
$ws.Range("P3").Value = 0.1005722428790014
$ws.Range("Q3").Value = 1.23881694534
$ws.Range("R3").Value = 11.14935250806
$ws.Range("S3").Value = 0.09768528951377062
$ws.Range("T3").Value = 0.1005722428790014

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.883656666666667
$ws.Range("H4").Value = 5.650970000000001
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 0.5797745
$ws.Range("N4").Value = 1.159549
$ws.Range("O4").Value = 0.08611580937010824
$ws.Range("P4").Value = 0.0591072285213179
$ws.Range("Q4").Value = 1.092096102088334
$ws.Range("R4").Value = 6.552576612530001
$ws.Range("S4").Value = 0.08611580937010824
$ws.Range("T4").Value = 0.0591072285213179
